$wb = $excel.ActiveWorkbook

# The "UK" sheet is a near-duplicate of the existing "Netherlands" sheet,
# so build it the same way a human would: copy the Netherlands tab to the
# end of the workbook, rename it, then touch up the handful of cells that
# differ.
$src = $wb.Worksheets.Item("Netherlands")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# New shared strings - set B4 first so "NGC-2741/T2176 " lands before
# "UK Market" in the shared string table (matches original authoring order).
$newSheet.Range("B4").Value = "NGC-2741/T2176 "
$newSheet.Range("B2").Value = "UK Market"

# Reorder accessory rows 10-13: MX-DPBX / MX-BBX now come before
# POS800-S / PR1D2-Unmonitored.
$newSheet.Range("A10").Value = "MX-DPBX"
$newSheet.Range("A11").Value = "MX-BBX"
$newSheet.Range("A12").Value = "POS800-S"
$newSheet.Range("A13").Value = "PR1D2-Unmonitored"

# Row 2 on Netherlands has an explicit wrapped-text height; the new sheet
# goes back to the sheet's default row height.
$newSheet.Rows.Item(2).AutoFit()

# Match the saved selection/active-cell state of the new tab.
$newSheet.Range("B4").Select()
